$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hyperlinks from A2 and B2
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()

# Clear A2 entirely (value + formatting)
$ws.Range("A2").Clear()

# Clear B2's contents but keep its style/formatting
$ws.Range("B2").ClearContents()

# Select B2 as the active cell
$ws.Range("B2").Select()
